# Natmi following Dr Hou advice
# Rebuild the LR-pairs result table (rows 2-13) as the full cross product of
# 4 sending clusters (ECs, FAPs, M2, sCs) x 3 target clusters (ECs, FAPs, sCs)
# for the Angpt2 -> Tie1 pair, replacing the original 4 rows (target = ECs only).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Angpt2","Tie1","ECs",3,1,5.709516000000001,17.128548,0.394755842864453,0.394755842864453,3,1,55.49161133333333,166.474834,0.9954988457613663,0.9954988457613663,316.830242773448,2851.472184961032,0.3929789859291182,0.3929789859291182),
    @("ECs","Angpt2","Tie1","FAPs",3,1,5.709516000000001,17.128548,0.394755842864453,0.394755842864453,1,0.3333333333333333,0.015201,0.045603,0.0002727002801111403,0.0002727002801111403,0.08679035271600001,0.7811131744440001,0.0001076500289246456,0.0001076500289246456),
    @("ECs","Angpt2","Tie1","sCs",3,1,5.709516000000001,17.128548,0.394755842864453,0.394755842864453,2,0.6666666666666666,0.2357046666666666,0.7071139999999999,0.00422845395852266,0.00422845395852266,1.345759565608,12.111836090472,0.001669206906410145,0.001669206906410145),
    @("FAPs","Angpt2","Tie1","ECs",3,1,3.281282666666666,9.843847999999999,0.2268678299100168,0.2268678299100168,3,1,55.49161133333333,166.474834,0.9954988457613663,0.9954988457613663,182.0836624134702,1638.752961721232,0.2258466628158077,0.2258466628158077),
    @("FAPs","Angpt2","Tie1","FAPs",3,1,3.281282666666666,9.843847999999999,0.2268678299100168,0.2268678299100168,1,0.3333333333333333,0.015201,0.045603,0.0002727002801111403,0.0002727002801111403,0.049878777816,0.4489090003439999,0.00006186692076466812,0.00006186692076466812),
    @("FAPs","Angpt2","Tie1","sCs",3,1,3.281282666666666,9.843847999999999,0.2268678299100168,0.2268678299100168,2,0.6666666666666666,0.2357046666666666,0.7071139999999999,0.00422845395852266,0.00422845395852266,0.7734136371857776,6.960722734671998,0.000959300173444456,0.0009593001734444561),
    @("M2","Angpt2","Tie1","ECs",3,1,1.388848333333333,4.166545,0.09602495105292475,0.09602495105292475,3,1,55.49161133333333,166.474834,0.9954988457613663,0.9954988457613663,77.06943191428111,693.6248872285299,0.09559272793747828,0.09559272793747828),
    @("M2","Angpt2","Tie1","FAPs",3,1,1.388848333333333,4.166545,0.09602495105292475,0.09602495105292475,1,0.3333333333333333,0.015201,0.045603,0.0002727002801111403,0.0002727002801111403,0.021111883515,0.190006951635,0.00002618603104979112,0.00002618603104979112),
    @("M2","Angpt2","Tie1","sCs",3,1,1.388848333333333,4.166545,0.09602495105292475,0.09602495105292475,2,0.6666666666666666,0.2357046666666666,0.7071139999999999,0.00422845395852266,0.00422845395852266,0.3273580334588889,2.94622230113,0.0004060370843966843,0.0004060370843966843),
    @("sCs","Angpt2","Tie1","ECs",3,1,4.083764,12.251292,0.2823513761726055,0.2823513761726055,3,1,55.49161133333333,166.474834,0.9954988457613663,0.9954988457613663,226.6146446650587,2039.531801985528,0.2810804690789621,0.2810804690789621),
    @("sCs","Angpt2","Tie1","FAPs",3,1,4.083764,12.251292,0.2823513761726055,0.2823513761726055,1,0.3333333333333333,0.015201,0.045603,0.0002727002801111403,0.0002727002801111403,0.062077296564,0.558695669076,0.00007699729937203546,0.00007699729937203546),
    @("sCs","Angpt2","Tie1","sCs",3,1,4.083764,12.251292,0.2823513761726055,0.2823513761726055,2,0.6666666666666666,0.2357046666666666,0.7071139999999999,0.00422845395852266,0.00422845395852266,0.9625622323653333,8.663060091287999,0.001193909794271374,0.001193909794271374)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
